# Insert a new record at row 58 of Sheet1 (pushing the existing rows
# 58:74 down to 59:75), then populate the new row with the weekly price
# entry for "Navel Late" fruit from "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58:74 down to 59:75, inserting a fresh (blank) row 58.
$ws.Rows.Item(58).Insert()

# Fill in the new row 58 with the new weekly record.
$ws.Cells.Item(58, 1).Value = 1
$ws.Cells.Item(58, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(58, 4).Value = 44524
$ws.Cells.Item(58, 5).Value = 15
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100102
$ws.Cells.Item(58, 8).Value = "Cítricos"
$ws.Cells.Item(58, 9).Value = 100102005
$ws.Cells.Item(58, 10).Value = "Naranja"
$ws.Cells.Item(58, 11).Value = "Navel Late"
$ws.Cells.Item(58, 12).Value = "Segunda"
$ws.Cells.Item(58, 13).Value = 300
$ws.Cells.Item(58, 14).Value = 700
$ws.Cells.Item(58, 15).Value = 750
$ws.Cells.Item(58, 16).Value = 725
$ws.Cells.Item(58, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(58, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(58, 19).Value = 725
$ws.Cells.Item(58, 20).Value = 1

# Keep the date cell formatted the same way as the rest of column D
# (the Insert() above already copies D57's format onto the new D58,
# but set it explicitly to be safe).
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
